# Refresh the cryptocurrency price/volume snapshot (rows 2-51) to the
# latest scraped values, per the "Updated cryptos list ... GitHub Actions"
# commit. A handful of rows also swap which coin occupies that rank
# (Algorand/TrustWalletToken and Aave/WEMIXTOKEN traded places), so those
# rows' Coin name + Link cells are updated alongside Price/Volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a target cell + its new text. Price-column values that look
# like plain numbers (e.g. "53.60", "1.009") are written with a leading
# apostrophe so Excel's COM layer stores them as literal text instead of
# silently coercing them to a Double and dropping the trailing zero / exact
# decimal formatting (e.g. 53.60 -> 53.6). The cell's Style is then reset to
# "Normal" so no stray quote-prefixed/text number format lingers on it.
$updates = @(
    @{ Cell = 'D2'; Value = '29.637.35' },
    @{ Cell = 'E2'; Value = '  -1.50%  ' },
    @{ Cell = 'D3'; Value = '2.092.48' },
    @{ Cell = 'E3'; Value = '  -0.50%  ' },
    @{ Cell = 'D4'; Value = '1.009' },
    @{ Cell = 'E4'; Value = '  +0.25%  ' },
    @{ Cell = 'D5'; Value = '343.21' },
    @{ Cell = 'E5'; Value = '  -1.26%  ' },
    @{ Cell = 'D6'; Value = '1.008' },
    @{ Cell = 'E6'; Value = '  +0.31%  ' },
    @{ Cell = 'D7'; Value = '0.5150' },
    @{ Cell = 'E7'; Value = '  -0.23%  ' },
    @{ Cell = 'D8'; Value = '0.4359' },
    @{ Cell = 'E8'; Value = '  -2.45%  ' },
    @{ Cell = 'D9'; Value = '53.60' },
    @{ Cell = 'E9'; Value = '  +1.99%  ' },
    @{ Cell = 'D10'; Value = '0.09153' },
    @{ Cell = 'E10'; Value = '  +2.43%  ' },
    @{ Cell = 'D11'; Value = '1.163' },
    @{ Cell = 'E11'; Value = '  -0.87%  ' },
    @{ Cell = 'D12'; Value = '24.50' },
    @{ Cell = 'E12'; Value = '  -4.13%  ' },
    @{ Cell = 'D13'; Value = '2.097.45' },
    @{ Cell = 'E13'; Value = '  -0.06%  ' },
    @{ Cell = 'D14'; Value = '6.739' },
    @{ Cell = 'E14'; Value = '  +0.19%  ' },
    @{ Cell = 'D15'; Value = '8.193' },
    @{ Cell = 'E15'; Value = '  +2.19%  ' },
    @{ Cell = 'D16'; Value = '100.20' },
    @{ Cell = 'E16'; Value = '  +0.54%  ' },
    @{ Cell = 'D17'; Value = '0.00001150' },
    @{ Cell = 'E17'; Value = '  +0.43%  ' },
    @{ Cell = 'D18'; Value = '1.009' },
    @{ Cell = 'E18'; Value = '  +0.23%  ' },
    @{ Cell = 'D19'; Value = '20.85' },
    @{ Cell = 'E19'; Value = '  +2.31%  ' },
    @{ Cell = 'D20'; Value = '0.06674' },
    @{ Cell = 'E20'; Value = '  -0.03%  ' },
    @{ Cell = 'E21'; Value = '  +0.32%  ' },
    @{ Cell = 'D22'; Value = '6.182' },
    @{ Cell = 'E22'; Value = '  -0.15%  ' },
    @{ Cell = 'D23'; Value = '29.696.40' },
    @{ Cell = 'E23'; Value = '  -1.59%  ' },
    @{ Cell = 'D24'; Value = '12.38' },
    @{ Cell = 'E24'; Value = '  -3.55%  ' },
    @{ Cell = 'D25'; Value = '2.303' },
    @{ Cell = 'E25'; Value = '  -2.16%  ' },
    @{ Cell = 'D26'; Value = '2.345.29' },
    @{ Cell = 'E26'; Value = '  -0.15%  ' },
    @{ Cell = 'D27'; Value = '21.80' },
    @{ Cell = 'E27'; Value = '  -0.74%  ' },
    @{ Cell = 'D28'; Value = '160.88' },
    @{ Cell = 'E28'; Value = '  -1.54%  ' },
    @{ Cell = 'D29'; Value = '2.480' },
    @{ Cell = 'E29'; Value = '  -2.32%  ' },
    @{ Cell = 'D30'; Value = '133.09' },
    @{ Cell = 'E30'; Value = '  -0.14%  ' },
    @{ Cell = 'D31'; Value = '1.122' },
    @{ Cell = 'E31'; Value = '  -4.98%  ' },
    @{ Cell = 'D32'; Value = '0.1047' },
    @{ Cell = 'E32'; Value = '  -1.82%  ' },
    @{ Cell = 'D33'; Value = '1.649' },
    @{ Cell = 'E33'; Value = '  +0.38%  ' },
    @{ Cell = 'D34'; Value = '6.175' },
    @{ Cell = 'E34'; Value = '  -1.10%  ' },
    @{ Cell = 'D35'; Value = '3.943' },
    @{ Cell = 'E35'; Value = '  -0.50%  ' },
    @{ Cell = 'D36'; Value = '6.313' },
    @{ Cell = 'E36'; Value = '  +6.12%  ' },
    @{ Cell = 'D37'; Value = '10.31' },
    @{ Cell = 'E37'; Value = '  +1.92%  ' },
    @{ Cell = 'D38'; Value = '0.02572' },
    @{ Cell = 'E38'; Value = '  -0.52%  ' },
    @{ Cell = 'D39'; Value = '0.06667' },
    @{ Cell = 'E39'; Value = '  -2.40%  ' },
    @{ Cell = 'D40'; Value = '0.6981' },
    @{ Cell = 'E40'; Value = '  +2.38%  ' },
    @{ Cell = 'D41'; Value = '12.42' },
    @{ Cell = 'E41'; Value = '  -0.66%  ' },
    @{ Cell = 'B42'; Value = 'TrustWalletToken' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Cell = 'D42'; Value = '1.319' },
    @{ Cell = 'E42'; Value = '  +5.27%  ' },
    @{ Cell = 'B43'; Value = 'Algorand' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D43'; Value = '0.2209' },
    @{ Cell = 'E43'; Value = '  -4.78%  ' },
    @{ Cell = 'D44'; Value = '0.6884' },
    @{ Cell = 'E44'; Value = '  +7.75%  ' },
    @{ Cell = 'D45'; Value = '14.34' },
    @{ Cell = 'E45'; Value = '  +0.45%  ' },
    @{ Cell = 'D46'; Value = '2.297' },
    @{ Cell = 'E46'; Value = '  +0.58%  ' },
    @{ Cell = 'D47'; Value = '3.620' },
    @{ Cell = 'E47'; Value = '  -1.12%  ' },
    @{ Cell = 'D48'; Value = '0.00000000344' },
    @{ Cell = 'E48'; Value = '  -4.78%  ' },
    @{ Cell = 'D49'; Value = '1.211' },
    @{ Cell = 'E49'; Value = '  -0.98%  ' },
    @{ Cell = 'B50'; Value = 'WEMIXTOKEN' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D50'; Value = '1.182' },
    @{ Cell = 'E50'; Value = '  +1.21%  ' },
    @{ Cell = 'B51'; Value = 'Aave' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Cell = 'D51'; Value = '1.182' },
    @{ Cell = 'E51'; Value = '  -2.42%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $value = $u.Value
    $looksNumeric = $value -match '^[+-]?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $cell.Value = "'" + $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}
